{"js": "// 1. Development skills bullet: add \"Spring Boot, \" before \"DropWizard.\"\n{\n  const results = context.document.body.search(\"Guava, Commons, DropWizard.\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"Guava, Commons, Spring Boot, DropWizard.\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 2. Split the Mendeley \"Part of multiple cross-functional teams...\" bullet into\n//    three separate bullet points (same list / style), dropping the \"** \" separators.\n{\n  const results = context.document.body.search(\n    \"Part of multiple cross-functional teams for Mendeley, through a major technology refresh. ** Acquisition and Onboarding: Migrating local Oauth2 sign-in to federated OpenID Connect solution. High-volume, mission-critical services. ** Building new services for flagship 'Reference Manager 2' product - mix of client-facing and message-processing.\",\n    { matchCase: true, matchWholeWord: false }\n  );\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    const firstPara = results.items[0].paragraphs.getFirst();\n    firstPara.load(\"style\");\n    await context.sync();\n    const listStyle = firstPara.style;\n\n    // Trim the original paragraph down to just the first sentence.\n    firstPara.insertText(\n      \"Part of multiple cross-functional teams for Mendeley, through a major technology refresh.\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n\n    // Insert the second bullet right after the (now shortened) first paragraph.\n    const secondPara = firstPara.insertParagraph(\n      \"Acquisition and Onboarding: Migrating local Oauth2 sign-in to federated OpenID Connect solution. High-volume, mission-critical services.\",\n      Word.InsertLocation.after\n    );\n    secondPara.style = listStyle;\n    await context.sync();\n    secondPara.attachToList(1006, 0);\n    await context.sync();\n\n    // Insert the third bullet right after the second.\n    const thirdPara = secondPara.insertParagraph(\n      \"Building new services for flagship 'Reference Manager 2' product - mix of client-facing and message-processing.\",\n      Word.InsertLocation.after\n    );\n    thirdPara.style = listStyle;\n    await context.sync();\n    thirdPara.attachToList(1006, 0);\n    await context.sync();\n  }\n}\n\n// 3. Technologies used/learned bullet: \"OpenID Connect\" -> \"OpenID Connect (OIDC)\",\n//    \"Payments (Adyen)\" -> \"Payments integration (Adyen)\".\n{\n  const results = context.document.body.search(\n    \"Technologies used/learned: Java 8, Dropwizard, Kibana, Redis, TDD, BDD, Oauth2, OpenID Connect, RxJava, Docker/ECS, Terraform, AWS, Payments (Adyen).\",\n    { matchCase: true, matchWholeWord: false }\n  );\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\n      \"Technologies used/learned: Java 8, Dropwizard, Kibana, Redis, TDD, BDD, Oauth2, OpenID Connect (OIDC), RxJava, Docker/ECS, Terraform, AWS, Payments integration (Adyen).\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n// 4. BBC bullet: \"Java/Camel/Jersey applications\" -> \"Java (Camel/Jersey/Spring Boot) applications\"\n{\n  const results = context.document.body.search(\"Java/Camel/Jersey applications\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"Java (Camel/Jersey/Spring Boot) applications\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 5. Remove the \"1993-1999 Granada Learning Ltd\" heading and its bullet entirely.\n{\n  const results = context.document.body.search(\"1993-1999 Granada Learning Ltd\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    const headingPara = results.items[0].paragraphs.getFirst();\n    const bulletPara = headingPara.getNext();\n    bulletPara.delete();\n    headingPara.delete();\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Development skills bullet: add \"Spring Boot, \" before \"DropWizard.\"\n$find1 = $d.Content.Find\n$find1.Execute(\n  \"Guava, Commons, DropWizard.\", $false, $false, $false, $false, $false, $true, 1, $false,\n  \"Guava, Commons, Spring Boot, DropWizard.\", 2\n)\n\n# 2. Split the Mendeley \"Part of multiple cross-functional teams...\" bullet into\n#    three separate bullet points (same list / style), dropping the \"** \" separators.\n$mendeleyPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  if ($d.Paragraphs.Item($i).Range.Text -like \"Part of multiple cross-functional teams for Mendeley*\") {\n    $mendeleyPara = $d.Paragraphs.Item($i)\n    break\n  }\n}\n\nif ($mendeleyPara -ne $null) {\n  # Insert two new (empty) paragraphs right after it, inheriting style/list formatting.\n  $mendeleyPara.Range.InsertParagraphAfter()\n  $secondPara = $mendeleyPara.Next()\n  $secondPara.Range.InsertParagraphAfter()\n  $thirdPara = $secondPara.Next()\n\n  $secondPara.Range.Text = \"Acquisition and Onboarding: Migrating local Oauth2 sign-in to federated OpenID Connect solution. High-volume, mission-critical services.\"\n  $thirdPara.Range.Text = \"Building new services for flagship 'Reference Manager 2' product - mix of client-facing and message-processing.\"\n\n  # Trim the original paragraph down to just the first sentence (exclude the paragraph mark).\n  $r = $mendeleyPara.Range\n  $r.End = $r.End - 1\n  $r.Text = \"Part of multiple cross-functional teams for Mendeley, through a major technology refresh.\"\n}\n\n# 3. Technologies used/learned bullet: \"OpenID Connect\" -> \"OpenID Connect (OIDC)\",\n#    \"Payments (Adyen)\" -> \"Payments integration (Adyen)\".\n$find2 = $d.Content.Find\n$find2.Execute(\n  \"OpenID Connect, RxJava\", $false, $false, $false, $false, $false, $true, 1, $false,\n  \"OpenID Connect (OIDC), RxJava\", 2\n)\n$find3 = $d.Content.Find\n$find3.Execute(\n  \"Payments (Adyen).\", $false, $false, $false, $false, $false, $true, 1, $false,\n  \"Payments integration (Adyen).\", 2\n)\n\n# 4. BBC bullet: \"Java/Camel/Jersey applications\" -> \"Java (Camel/Jersey/Spring Boot) applications\"\n$find4 = $d.Content.Find\n$find4.Execute(\n  \"Java/Camel/Jersey applications\", $false, $false, $false, $false, $false, $true, 1, $false,\n  \"Java (Camel/Jersey/Spring Boot) applications\", 2\n)\n\n# 5. Remove the \"1993-1999 Granada Learning Ltd\" heading and its bullet entirely.\n$granadaHeading = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  if ($d.Paragraphs.Item($i).Range.Text -like \"1993-1999 Granada Learning Ltd*\") {\n    $granadaHeading = $d.Paragraphs.Item($i)\n    break\n  }\n}\n\nif ($granadaHeading -ne $null) {\n  $granadaBullet = $granadaHeading.Next()\n  $granadaBullet.Range.Delete()\n  $granadaHeading.Range.Delete()\n}\n"}
